$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column N (14); existing N..Q shift right to O..R
$ws.Columns.Item(14).Insert()

# Give the newly inserted header cell (N1) its text, matching the style of
# its neighboring header cells (gray header style already used by row 1)
$ws.Cells.Item(1, 14).Value = "利率種類"

# Narrow column M (13) and size the new column N (14)
$ws.Columns.Item(13).ColumnWidth = 10.3
$ws.Columns.Item(14).ColumnWidth = 10.85

# Update the active selection shown in the sheet view
$null = $ws.Range("M7").Select()

# Update the hidden _FilterDatabase defined name so it covers the new range
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=正常件!`$A`$1:`$Q`$1"
    }
}
